$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header for column A
$ws.Range("A1").Value = "location"

# Delete the last row (row 18), shrinking the used range to A1:B17
$ws.Rows.Item(18).Delete() | Out-Null

# Replace column B (activity name) values with the new English names,
# row-by-row (A column grid codes are unchanged).
$names = @(
    "activity_name",
    "a",
    "b",
    "c",
    "d",
    "e",
    "a",
    "b",
    "d",
    "h",
    "c",
    "d",
    "f",
    "b",
    "d",
    "e",
    "j"
)

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 2).Value = $names[$i]
}

# Reset selection to default (A1)
$ws.Range("A1").Select() | Out-Null
